$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Move the "Changes Made" block up from rows 15-17 to rows 13-15
#    (frees up rows 16-17 and keeps row 15 for the last line)
# ------------------------------------------------------------------
$changesMade   = $ws.Range("C15").Value2
$addedDesc     = $ws.Range("C16").Value2
$addedCoords   = $ws.Range("C17").Value2

$ws.Range("C13").Value = $changesMade
$ws.Range("C14").Value = $addedDesc
$ws.Range("C15").Value = $addedCoords
$ws.Range("C16").ClearContents()
$ws.Range("C17").ClearContents()

# ------------------------------------------------------------------
# 2. Highlight the legend row (row 2) with a yellow fill
# ------------------------------------------------------------------
$ws.Range("A2:C2").Interior.Color = 65535

# ------------------------------------------------------------------
# 3. Add a second legend row (row 18) above the new "Table13" table,
#    duplicating row 2's content & new yellow fill style
# ------------------------------------------------------------------
$ws.Range("A18").Value = $ws.Range("A2").Value2
$ws.Range("B18").Value = $ws.Range("B2").Value2
$ws.Range("C18").Value = $ws.Range("C2").Value2
$ws.Range("A18:C18").Interior.Color = 65535

# ------------------------------------------------------------------
# 4. New survey block #2 -> rows 19-25 (becomes Table13)
# ------------------------------------------------------------------
$ws.Range("A19").Value = "Category:"
$ws.Range("B19").Value = "Rating (1-10):"
$ws.Range("C19").Value = "Comments:"

$ws.Range("A20").Value = "Design"
$ws.Range("A21").Value = "Instruction clarity"
$ws.Range("A22").Value = "Ease of use"
$ws.Range("A23").Value = "Cool factor"

$ws.Range("A24").Value = "Total"
$ws.Range("B24").Formula = "=SUM(B20:B23)"

$ws.Range("A25").Value = "Max possible score"
$ws.Range("B25").Value = 40

# ------------------------------------------------------------------
# 5. New survey block #3 -> rows 27-33 (becomes Table134)
# ------------------------------------------------------------------
$ws.Range("A27").Value = "Category:"
$ws.Range("B27").Value = "Rating (1-10):"
$ws.Range("C27").Value = "Comments:"

$ws.Range("A28").Value = "Design"
$ws.Range("A29").Value = "Instruction clarity"
$ws.Range("A30").Value = "Ease of use"
$ws.Range("A31").Value = "Cool factor"

$ws.Range("A32").Value = "Total"
$ws.Range("B32").Formula = "=SUM(B28:B31)"

$ws.Range("A33").Value = "Max possible score"
$ws.Range("B33").Value = 40

# ------------------------------------------------------------------
# 6. Turn the two new ranges into real Excel Tables (ListObjects),
#    matching the style used by the original Table1
# ------------------------------------------------------------------
$tbl2 = $ws.ListObjects.Add(1, $ws.Range("A19:C25"), $null, 1)
$tbl2.Name = "Table13"
$tbl2.TableStyle = "TableStyleLight8"

$tbl3 = $ws.ListObjects.Add(1, $ws.Range("A27:C33"), $null, 1)
$tbl3.Name = "Table134"
$tbl3.TableStyle = "TableStyleLight8"

# ------------------------------------------------------------------
# 7. View state - zoom in to 130% and leave the selection on B20,
#    matching where the author was working
# ------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 130
$ws.Range("B20").Select() | Out-Null
